$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J, matching the style of the
# existing header cells (bold, bordered, centered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-4, columns I and J
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 3
